$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4: fill in B4:F4 (A4 already has value 3)
$ws.Range("B4").Value = 44316
$ws.Range("C4").Value = "Julio Edwin Mora Ardila"
$ws.Range("D4").Value = "Creación repositorio"
$ws.Range("E4").Value = "Diagrama de base de datos"
$ws.Range("F4").Value = "Ninguna"

# Row 5: new row
$ws.Range("A5").Value = 4
$ws.Range("B5").Value = 44318
$ws.Range("C5").Value = "Juan Carlos Rojas Buitrago"
$ws.Range("D5").Value = "Alimentación repositorio"
$ws.Range("E5").Value = "Alimentación HU en Jira"
$ws.Range("F5").Value = "Ninguna"

$ws.Columns.Item(4).ColumnWidth = 17
$ws.Columns.Item(5).ColumnWidth = 26.3

$ws.Range("F5").Select()
